$wb = $excel.ActiveWorkbook

# --- sheet2 ("IMG_5297.JPG"): clear the stray empty F2/F3 stub cells ---
$ws2 = $wb.Worksheets.Item("IMG_5297.JPG")
$ws2.Range("F2").ClearContents()
$ws2.Range("F3").ClearContents()

# --- helper to stamp the standard header row on a freshly added sheet ---
function Set-Header($ws) {
    $ws.Range("A1").Value = "project_id"
    $ws.Range("B1").Value = "image_file_name"
    $ws.Range("C1").Value = "colony_label"
    $ws.Range("D1").Value = "number_of_colonies"
    $ws.Range("E1").Value = "average_size"
    $ws.Range("F1").Value = "std_dev_size"

    $hdr = $ws.Range("A1:F1")
    $hdr.Font.Bold = $true
    $hdr.HorizontalAlignment = -4108
    $hdr.VerticalAlignment = -4160
    $hdr.Borders.LineStyle = 1
    $hdr.Borders.Weight = 2
}

# --- new sheet 3: "IMG_5265.jpg" ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "IMG_5265.jpg"
Set-Header $ws3

$ws3.Range("A2").Value = 1
$ws3.Range("B2").Value = "IMG_5264.jpg"
$ws3.Range("C2").Value = "Kbd17-1"
$ws3.Range("D2").Value = 115
$ws3.Range("E2").Value = 944.1130000000001
$ws3.Range("F2").Value = 463.726

$ws3.Range("A3").Value = 2
$ws3.Range("B3").Value = "IMG_5265.jpg"
$ws3.Range("C3").Value = "Kbd17-2"
$ws3.Range("D3").Value = 19
$ws3.Range("E3").Value = 685.158
$ws3.Range("F3").Value = 464.293

# --- new sheet 4: "IMG_5297.JPG1" ---
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add($null, $lastSheet2)
$ws4.Name = "IMG_5297.JPG1"
Set-Header $ws4

$ws4.Range("A2").Value = 1
$ws4.Range("B2").Value = "IMG_5264.jpg"
$ws4.Range("C2").Value = "Kbd17-1"
$ws4.Range("D2").Value = 115
$ws4.Range("E2").Value = 944.1130000000001
$ws4.Range("F2").Value = 463.726

$ws4.Range("A3").Value = 2
$ws4.Range("B3").Value = "IMG_5265.jpg"
$ws4.Range("C3").Value = "Kbd17-2"
$ws4.Range("D3").Value = 19
$ws4.Range("E3").Value = 685.158
$ws4.Range("F3").Value = 464.293

$ws4.Range("A4").Value = 3
$ws4.Range("B4").Value = "IMG_5297.JPG"
$ws4.Range("C4").Value = "Kbd5-2"
$ws4.Range("D4").Value = 1277
$ws4.Range("E4").Value = 276.432
$ws4.Range("F4").Value = 189.652
